$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New label text for column A (rows 2-10), replacing the SN1/SN2/SN3 naming
# scheme with FA1/FA2/FA3, and collapsing the 12 distinct labels down to the
# 9 labels actually used now (dropping the "[M-(SNx-H2O)+H]+" group entirely).
# Values are written grouped by FA-number (as if via a series of Find/Replace
# passes: SN1->FA1, then SN2->FA2, then SN3->FA3) so that the resulting
# shared-string table ordering matches the original edit.
$ws.Range("A2").Value = "FA1_[FA-H2O+H]+"
$ws.Range("A5").Value = "[MG(FA1)-H2O+H]+"
$ws.Range("A8").Value = "[M-(FA1)+H]+"

$ws.Range("A3").Value = "FA2_[FA-H2O+H]+"
$ws.Range("A6").Value = "[MG(FA2)-H2O+H]+"
$ws.Range("A9").Value = "[M-(FA2)+H]+"

$ws.Range("A4").Value = "FA3_[FA-H2O+H]+"
$ws.Range("A7").Value = "[MG(FA3)-H2O+H]+"
$ws.Range("A10").Value = "[M-(FA3)+H]+"

# Update the weight column (B) for the first two groups of rows: the old
# value of 1 is corrected to 1.6 to fix the [M+Na]+ / [M+H]+ weighting bug.
$ws.Range("B2").Value = 1.6
$ws.Range("B3").Value = 1.6
$ws.Range("B4").Value = 1.6
$ws.Range("B5").Value = 1.6
$ws.Range("B6").Value = 1.6
$ws.Range("B7").Value = 1.6

# Remove the now-unused rows 11-13 (the old "[M-(SNx-H2O)+H]+" group), which
# shifts nothing else since they were already the last rows.
$ws.Range("A11:C13").Delete()

# Update the selected cell to match the saved view state.
$ws.Range("A10").Select()
